$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, Open(D), Close(E), High(F), Low(G), SharesOutstanding(H)
$rows = @(
    @(2, 47.25, 39.38999938964844, 52.5099983215332, 38.75, 63060655),
    @(3, 43.36000061035156, 26.39999961853028, 43.4900016784668, 23.65999984741211, 63060655),
    @(4, 21.64999961853028, 22.25, 25.02000045776367, 20.60000038146973, 63060655),
    @(5, 28.39999961853028, 20.95000076293945, 28.54999923706055, 18.94000053405762, 63060655),
    @(6, 19.78000068664551, 21, 21.88999938964844, 19.20999908447266, 63060655),
    @(7, 30.67000007629395, 32.16999816894531, 32.40999984741211, 28.68000030517578, 63060655),
    @(8, 41.54999923706055, 32.65999984741211, 43.36000061035156, 32.09000015258789, 63060655),
    @(9, 38.5, 41.77999877929688, 43.02999877929688, 37.56000137329102, 63060655),
    @(10, 32.90000152587891, 35.40999984741211, 36.15999984741211, 31.61000061035156, 63060655),
    @(11, 30.01000022888184, 32.52999877929688, 33.31999969482422, 29.29999923706055, 63060655),
    @(12, 43.79999923706055, 46.72000122070312, 47.20000076293945, 42.91999816894531, 63060655),
    @(13, 42.06000137329102, 43.81999969482422, 46.13999938964844, 42.02000045776367, 63060655),
    @(14, 41.59000015258789, 44.84999847412109, 45.81999969482422, 40.16999816894531, 63060655),
    @(15, 38.77999877929688, 36.88000106811523, 41.41999816894531, 36.40999984741211, 63060655),
    @(16, 49.65000152587891, 42.81999969482422, 49.93000030517578, 39.5099983215332, 63060655),
    @(17, 34.18999862670898, 36.41999816894531, 37.29999923706055, 32.84000015258789, 63060655),
    @(18, 34.90000152587891, 40.06000137329102, 40.43999862670898, 34.61000061035156, 63060655),
    @(19, 34.75, 35.04999923706055, 36.31999969482422, 33.77000045776367, 63060655),
    @(20, 34.95000076293945, 34.5099983215332, 35.61999893188477, 32.16999816894531, 63060655),
    @(21, 35, 32.59999847412109, 37.4900016784668, 32.56999969482422, 63060655),
    @(22, 17.01000022888184, 22.35000038146973, 24.53000068664551, 16.15999984741211, 63060655),
    @(23, 23.1200008392334, 24.97999954223633, 25.30999946594238, 21.64999961853028, 63060655),
    @(24, 20.18000030517578, 19.67000007629395, 21.89999961853028, 18.93000030517578, 63060655),
    @(25, 32.5099983215332, 32.59000015258789, 36.72999954223633, 29.21999931335449, 63060655),
    @(26, 39.63000106811523, 39.29999923706055, 42.45000076293945, 37.58000183105469, 63060655),
    @(27, 40.13999938964844, 37.40000152587891, 41.06999969482422, 35.29999923706055, 63060655),
    @(28, 37.43000030517578, 38.63000106811523, 40.79999923706055, 37.43000030517578, 63060655),
    @(29, 36.40000152587891, 34.54000091552734, 37.68000030517578, 31.60000038146973, 63060655),
    @(30, 34.18999862670898, 32.52999877929688, 35.70999908447266, 32.11000061035156, 63060655),
    @(31, 27.61000061035156, 30.65999984741211, 31.63999938964844, 26.28000068664551, 63060655),
    @(32, 34.08000183105469, 38.40999984741211, 39.27000045776367, 33.34000015258789, 63060655),
    @(33, 27.82999992370605, 31.51000022888184, 31.54000091552734, 27.40999984741211, 63060655),
    @(34, 30.38999938964844, 29.92000007629395, 31.22999954223633, 28.45000076293945, 63060655),
    @(35, 36.11000061035156, 45.04999923706055, 45.52999877929688, 35.61000061035156, 63060655),
    @(36, 41.63000106811523, 42.18999862670898, 44.41999816894531, 39.70999908447266, 63060655),
    @(37, 46.83000183105469, 43.72999954223633, 46.97999954223633, 43.04000091552734, 63060655),
    @(38, 39.40999984741211, 40.2400016784668, 41.58000183105469, 38.93000030517578, 63060655),
    @(39, 36.93999862670898, 36.43000030517578, 38.11999893188477, 34.7599983215332, 63060655),
    @(40, 35.13000106811523, 34.13999938964844, 35.43000030517578, 33.20999908447266, 63060655),
    @(41, 39.2400016784668, 39.93999862670898, 41.72000122070312, 37.95999908447266, 63060655),
    @(42, 36.72000122070312, 35.08000183105469, 38.54000091552734, 32.27999877929688, 63060655),
    @(43, 34.18999862670898, 34.43000030517578, 35.9900016784668, 33.95999908447266, 63060655)
)

foreach ($entry in $rows) {
    $r = $entry[0]
    $ws.Cells.Item($r, 4).Value = $entry[1]   # D: open_price
    $ws.Cells.Item($r, 5).Value = $entry[2]   # E: close_price
    $ws.Cells.Item($r, 6).Value = $entry[3]   # F: high_price
    $ws.Cells.Item($r, 7).Value = $entry[4]   # G: low_price
    $ws.Cells.Item($r, 8).Value = $entry[5]   # H: shares_outstanding
    $ws.Cells.Item($r, 9).Value = "YELP"      # I: fixed_ticker
}
